$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Clear the "Title" value (B5) — the "Forme pharmaceutique" text is removed.
$ws.Range("B5").Value = ""

# Update the "Date" value (B8) to reflect the new export timestamp.
$ws.Range("B8").Value = "2024-06-04T08:55:54+00:00"
